$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Стачки 204")

$ws.Range("B4:D4").NumberFormat = "@"

$ws.Range("A4").Value = "Стачки 204, кв. 30"
$ws.Range("B4").Value = "453"
$ws.Range("C4").Value = "68"
$ws.Range("D4").Value = "19473"
$ws.Range("E4").Value = "08.11.2023 в 14:58:41"
$ws.Range("F4").Value = "Alex Pol ID 128446192"
